# Weekly price-sheet update: a new current-week record is published for
# "Feria Lagunitas de Puerto Montt - Pepino ensalada" and inserted at row
# 336, pushing the previously existing rows 336-373 down to 337-374
# (dimension grows from A1:R373 to A1:R374).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 336, shifting rows 336:373
# (and everything below) down by one - mirrors Excel's native
# "Insert Sheet Rows" behaviour.
$ws.Rows.Item(336).Insert()

# Populate the newly inserted row 336 with the new week's data. Columns
# that are identical to the record this new row is based on (A, B, C, E,
# F, G, H, I, J, N, O, Q, R) are filled in as well since the row starts
# out empty after the insert.
$ws.Range("A336").Value = 4
$ws.Range("B336").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C336").Value = 'Los Lagos'
$ws.Range("D336").Value = 44946
$ws.Range("E336").Value = 10
$ws.Range("F336").Value = 100112043
$ws.Range("G336").Value = 'Pepino ensalada'
$ws.Range("H336").Value = 'Sin especificar'
$ws.Range("I336").Value = 'Primera'
$ws.Range("J336").Value = 400
$ws.Range("K336").Value = 16000
$ws.Range("L336").Value = 18000
$ws.Range("M336").Value = 17000
$ws.Range("N336").Value = '$/caja 60 unidades'
$ws.Range("O336").Value = 'Región de Arica y Parinacota'
$ws.Range("P336").Value = 283
$ws.Range("Q336").Value = 60
$ws.Range("R336").Value = 'Hortaliza'
